$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new "Save" header column in H1, copying the header formatting
# (style) used by the other header cells (e.g. G1 "sum").
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the corresponding data value for the new "Save" column in H2.
$ws.Range("H2").Value = 0
